{"js": "// Office.js (Word JavaScript API) script.\n// Applies the resume wording tweaks described by the diff:\n//   1. Split the STEM/ASU bullet into two sentences.\n//   2. Rewrite the Go Fish bullet (single-player vs. two-player wording).\n//   3. Rewrite the Randoma11y bullet (\"their current tab\" / \"Displays\").\n//   4. Reorder the languages/technologies list (SQL moves before Jest, Cypress).\n//   5. Merge \"Visual \" + \"Studio, Git, Eclipse, XCode\" into one run that\n//      reads \"Visual Studio, Git, Eclipse\" (XCode dropped).\n\nconst body = context.document.body;\n\nasync function replaceOnce(searchText, replacementText) {\n  const results = body.search(searchText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n  results.items[0].insertText(replacementText, \"Replace\");\n  await context.sync();\n}\n\n// 1. STEM recruitment / ASU bullet.\nawait replaceOnce(\n  \"Documented trends in STEM recruitment and retention at ASU and created a research paper from our findings using LaTeX.\",\n  \"Documented trends in STEM recruitment and retention at ASU. Created a research paper from our findings using LaTeX.\"\n);\n\n// 2. Go Fish bullet (note the leading space \u2014 it follows the bold \"Go Fish:\" run).\nawait replaceOnce(\n  \" Developed a Go Fish game in C++ where two players can play against each other or play against the AI. The results for each game are written to the file system for keeping track of match history.\",\n  \" Developed a Go Fish game with C++ for one player to play against the computer. Results for each game are written to the file system to keep track of match history.\"\n);\n\n// 3. Randoma11y Chrome extension bullet.\nawait replaceOnce(\n  \"Created an extension that allows users to change the appearance of a tab in Chrome by toggling a randoma11y theme. Shows contrast ratio and WCAG AA/AAA information.\",\n  \"Created an extension that allows users to change the appearance of their current tab in Chrome by toggling a randoma11y theme. Displays contrast ratio and WCAG AA/AAA information.\"\n);\n\n// 4. Languages & Technologies list \u2014 move SQL earlier, right after Python.\nawait replaceOnce(\n  \"C++, Java, HTML, CSS, Sass, JavaScript, Node, React, Python, Jest, Cypress, SQL\",\n  \"C++, Java, HTML, CSS, Sass, JavaScript, Node, React, Python, SQL, Jest, Cypress\"\n);\n\n// 5. \"Visual \" + \"Studio, Git, Eclipse, XCode\" -> single run \"Visual Studio, Git, Eclipse\".\nawait replaceOnce(\n  \"Visual Studio, Git, Eclipse, XCode\",\n  \"Visual Studio, Git, Eclipse\"\n);\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Applies the resume wording tweaks described by the diff:\n#   1. Split the STEM/ASU bullet into two sentences.\n#   2. Rewrite the Go Fish bullet (single-player vs. two-player wording).\n#   3. Rewrite the Randoma11y bullet (\"their current tab\" / \"Displays\").\n#   4. Reorder the languages/technologies list (SQL moves before Jest, Cypress).\n#   5. Merge \"Visual \" + \"Studio, Git, Eclipse, XCode\" into one run that\n#      reads \"Visual Studio, Git, Eclipse\" (XCode dropped).\n\n$d = $word.ActiveDocument\n\nfunction Replace-Once($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n    if (-not $found) {\n        throw \"Text not found: $findText\"\n    }\n}\n\n# 1. STEM recruitment / ASU bullet.\nReplace-Once \"Documented trends in STEM recruitment and retention at ASU and created a research paper from our findings using LaTeX.\" \"Documented trends in STEM recruitment and retention at ASU. Created a research paper from our findings using LaTeX.\"\n\n# 2. Go Fish bullet (note the leading space -- it follows the bold \"Go Fish:\" run).\nReplace-Once \" Developed a Go Fish game in C++ where two players can play against each other or play against the AI. The results for each game are written to the file system for keeping track of match history.\" \" Developed a Go Fish game with C++ for one player to play against the computer. Results for each game are written to the file system to keep track of match history.\"\n\n# 3. Randoma11y Chrome extension bullet.\nReplace-Once \"Created an extension that allows users to change the appearance of a tab in Chrome by toggling a randoma11y theme. Shows contrast ratio and WCAG AA/AAA information.\" \"Created an extension that allows users to change the appearance of their current tab in Chrome by toggling a randoma11y theme. Displays contrast ratio and WCAG AA/AAA information.\"\n\n# 4. Languages & Technologies list -- move SQL earlier, right after Python.\nReplace-Once \"C++, Java, HTML, CSS, Sass, JavaScript, Node, React, Python, Jest, Cypress, SQL\" \"C++, Java, HTML, CSS, Sass, JavaScript, Node, React, Python, SQL, Jest, Cypress\"\n\n# 5. \"Visual \" + \"Studio, Git, Eclipse, XCode\" -> single run \"Visual Studio, Git, Eclipse\".\nReplace-Once \"Visual Studio, Git, Eclipse, XCode\" \"Visual Studio, Git, Eclipse\"\n"}
